$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Steps")

# Rename "human detection" -> "humanDetection" for the first guest's detect-human step
$ws.Range("G6").Value = "humanDetection"

# Fix "who " (trailing space) -> "who" for the first guest's Ask age step
$ws.Range("H11").Value = "who"

# Rename "living room" -> "livingRoom" (goTo living room step, guest 1)
$ws.Range("I14").Value = "livingRoom"

# Add a new location/livingRoom key-value pair to the "Ask to follow" step for guest 1
$ws.Range("J13").Value = "location"
$ws.Range("K13").Value = "livingRoom"

# Rename "door" -> "entrance" (goTo door step, guest 1)
$ws.Range("I26").Value = "entrance"

# Rename "living room" -> "livingRoom" (goTo living room step, guest 2)
$ws.Range("I39").Value = "livingRoom"

# Add a new location/livingRoom key-value pair to the "Ask to follow" step for guest 2
$ws.Range("J38").Value = "location"
$ws.Range("K38").Value = "livingRoom"

# Update view state: scroll position and current selection
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("J38:K38").Select()
